$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1727495209494947
$ws.Range("D2").Value = 0.127157671452153
$ws.Range("E2").Value = 0.1364255611091814
$ws.Range("F2").Value = 1.907882314879409
$ws.Range("G2").Value = 0.002513077859497646
$ws.Range("I2").Value = 1.345494556904121
$ws.Range("J2").Value = 0.1822171015749401
$ws.Range("K2").Value = 1.578902789843767
$ws.Range("L2").Value = 0.195507159602613
$ws.Range("M2").Value = 0.4910369676017794
$ws.Range("O2").Value = 5.022225227143252
$ws.Range("C3").Value = 0.1705136274307222
$ws.Range("D3").Value = 0.1248254744461406
$ws.Range("E3").Value = 0.1366810173988124
$ws.Range("F3").Value = 1.925036952250629
$ws.Range("G3").Value = 0.002516005637306301
$ws.Range("I3").Value = 1.358196868634025
$ws.Range("J3").Value = 0.1837903693103673
$ws.Range("K3").Value = 1.432124680810546
$ws.Range("L3").Value = 0.1967932989946952
$ws.Range("M3").Value = 0.4630080346250054
$ws.Range("O3").Value = 5.079078681962429
$ws.Range("C4").Value = 0.1691938727517339
$ws.Range("D4").Value = 0.1234219996810069
$ws.Range("E4").Value = 0.1368813955479578
$ws.Range("F4").Value = 1.936636649290755
$ws.Range("G4").Value = 0.002517898928674033
$ws.Range("I4").Value = 1.366761295343657
$ws.Range("J4").Value = 0.1848209741551354
$ws.Range("K4").Value = 1.341896932497178
$ws.Range("L4").Value = 0.1976431853464362
$ws.Range("M4").Value = 0.4458509285073688
$ws.Range("O4").Value = 5.116976913954886
$ws.Range("C5").Value = 0.1686694838806204
$ws.Range("D5").Value = 0.1228573039907204
$ws.Range("E5").Value = 0.136974022186104
$ws.Range("F5").Value = 1.941631812381466
$ws.Range("G5").Value = 0.002518694576528112
$ws.Range("I5").Value = 1.370443731627088
$ws.Range("J5").Value = 0.1852572164660238
$ws.Range("K5").Value = 1.305104387904322
$ws.Range("L5").Value = 0.1980046789056917
$ws.Range("M5").Value = 0.438873052026743
$ws.Range("O5").Value = 5.133172211641678
$ws.Range("C6").Value = 0.1685832222572969
$ws.Range("D6").Value = 0.1227639753736867
$ws.Range("E6").Value = 0.1369900661422339
$ws.Range("F6").Value = 1.942477452600002
$ws.Range("G6").Value = 0.002518828151976029
$ws.Range("I6").Value = 1.371066814385614
$ws.Range("J6").Value = 0.1853306368310701
$ws.Range("K6").Value = 1.298993638086898
$ws.Range("L6").Value = 0.1980656207958411
$ws.Range("M6").Value = 0.4377152286487842
$ws.Range("O6").Value = 5.135906803276072
$ws.Range("C7").Value = 0.1691867462208023
$ws.Range("D7").Value = 0.1234143546304409
$ws.Range("E7").Value = 0.1368826002881303
$ws.Range("F7").Value = 1.936702929953142
$ws.Range("G7").Value = 0.002517909561069287
$ws.Range("I7").Value = 1.366810179154783
$ws.Range("J7").Value = 0.1848267916049124
$ws.Range("K7").Value = 1.341400829238552
$ws.Range("L7").Value = 0.1976479991711226
$ws.Range("M7").Value = 0.4457567659316979
$ws.Range("O7").Value = 5.117192287301549
$ws.Range("C8").Value = 0.1719676070513714
$ws.Range("D8").Value = 0.1263476515080129
$ws.Range("E8").Value = 0.1365046236584941
$ws.Range("F8").Value = 1.913575824247175
$ws.Range("G8").Value = 0.002514067553667122
$ws.Range("I8").Value = 1.349715480725614
$ws.Range("J8").Value = 0.1827461632968816
$ws.Range("K8").Value = 1.52831701906922
$ws.Range("L8").Value = 0.1959381411335368
$ws.Range("M8").Value = 0.4813619364715294
$ws.Range("O8").Value = 5.041207493038883
$ws.Range("C9").Value = 0.1778394155848417
$ws.Range("D9").Value = 0.1323235460520209
$ws.Range("E9").Value = 0.1361077222680951
$ws.Range("F9").Value = 1.876689206823926
$ws.Range("G9").Value = 0.002507288861668879
$ws.Range("I9").Value = 1.322266138603474
$ws.Range("J9").Value = 0.1791780189937349
$ws.Range("K9").Value = 1.89393176442087
$ws.Range("L9").Value = 0.1930617131239956
$ws.Range("M9").Value = 0.5515834510862305
$ws.Range("O9").Value = 4.915943408352319
$ws.Range("C10").Value = 0.1824055377291103
$ws.Range("D10").Value = 0.1368476523574742
$ws.Range("E10").Value = 0.136024704541601
$ws.Range("F10").Value = 1.8547527812878
$ws.Range("G10").Value = 0.0025027645376865
$ws.Range("I10").Value = 1.30580489774708
$ws.Range("J10").Value = 0.1768676042910968
$ws.Range("K10").Value = 2.161888704201999
$ws.Range("L10").Value = 0.1912375902015171
$ws.Range("M10").Value = 0.6033991426643865
$ws.Range("O10").Value = 4.83840998310248
$ws.Range("C11").Value = 0.1845369223965321
$ws.Range("D11").Value = 0.1389342529038373
$ws.Range("E11").Value = 0.1360319600958348
$ws.Range("F11").Value = 1.845895472037014
$ws.Range("G11").Value = 0.00250080434306891
$ws.Range("I11").Value = 1.29912169261128
$ws.Range("J11").Value = 0.1758838764557424
$ws.Range("K11").Value = 2.283627819176957
$ws.Range("L11").Value = 0.1904702576986068
$ws.Range("M11").Value = 0.6270162665483099
$ws.Range("O11").Value = 4.806291476845786
$ws.Range("C12").Value = 0.1853517584823976
$ws.Range("D12").Value = 0.1397284451104213
$ws.Range("E12").Value = 0.1360411575807277
$ws.Range("F12").Value = 1.842702802131228
$ws.Range("G12").Value = 0.00250007607899249
$ws.Range("I12").Value = 1.296706774046676
$ws.Range("J12").Value = 0.1755210257471695
$ws.Range("K12").Value = 2.329702768439063
$ws.Range("L12").Value = 0.1901886500082455
$ws.Range("M12").Value = 0.6359656274428289
$ws.Range("O12").Value = 4.794582746043602
$ws.Range("C13").Value = 0.1851759264292951
$ws.Range("D13").Value = 0.1395572228811517
$ws.Range("E13").Value = 0.13603889025444
$ws.Range("F13").Value = 1.843383220713008
$ws.Range("G13").Value = 0.002500232301069938
$ws.Range("I13").Value = 1.297221715290185
$ws.Range("J13").Value = 0.1755987424261036
$ws.Range("K13").Value = 2.319780857273088
$ws.Range("L13").Value = 0.1902489008614907
$ws.Range("M13").Value = 0.6340379609748794
$ws.Range("O13").Value = 4.797084235736605
$ws.Range("C14").Value = 0.1846038049999237
$ws.Range("D14").Value = 0.1389995109418294
$ws.Range("E14").Value = 0.1360325876484154
$ws.Range("F14").Value = 1.845629573851625
$ws.Range("G14").Value = 0.002500744147708474
$ws.Range("I14").Value = 1.298920693137873
$ws.Range("J14").Value = 0.1758538308562301
$ws.Range("K14").Value = 2.287418949061305
$ws.Range("L14").Value = 0.1904469101092765
$ws.Range("M14").Value = 0.6277524172645457
$ws.Range("O14").Value = 4.805319093130436
$ws.Range("C15").Value = 0.1842543682655133
$ws.Range("D15").Value = 0.1386584207157426
$ws.Range("E15").Value = 0.1360295663866182
$ws.Range("F15").Value = 1.847026552263159
$ws.Range("G15").Value = 0.002501059492235045
$ws.Range("I15").Value = 1.299976458116369
$ws.Range("J15").Value = 0.1760113383242263
$ws.Range("K15").Value = 2.267593011839381
$ws.Range("L15").Value = 0.1905693633786925
$ws.Range("M15").Value = 0.6239031148408856
$ws.Range("O15").Value = 4.8104223041654
$ws.Range("C16").Value = 0.1822673318969663
$ws.Range("D16").Value = 0.1367118572505319
$ws.Range("E16").Value = 0.1360251342556076
$ws.Range("F16").Value = 1.855354207312246
$ws.Range("G16").Value = 0.002502894606334552
$ws.Range("I16").Value = 1.306257869583781
$ws.Range("J16").Value = 0.1769332463725934
$ws.Range("K16").Value = 2.153929407819135
$ws.Range("L16").Value = 0.1912889926471522
$ws.Range("M16").Value = 0.6018565861859457
$ws.Range("O16").Value = 4.840572479086319
$ws.Range("C17").Value = 0.1810621886688182
$ws.Range("D17").Value = 0.1355249723929148
$ws.Range("E17").Value = 0.1360339290940615
$ws.Range("F17").Value = 1.860750317873737
$ws.Range("G17").Value = 0.002504045428398561
$ws.Range("I17").Value = 1.31031758371342
$ws.Range("J17").Value = 0.1775160343605435
$ws.Range("K17").Value = 2.084158694291375
$ws.Range("L17").Value = 0.1917464486909317
$ws.Range("M17").Value = 0.5883431321011869
$ws.Range("O17").Value = 4.859876337688632
$ws.Range("C18").Value = 0.180374131065804
$ws.Range("D18").Value = 0.1348449998718309
$ws.Range("E18").Value = 0.1360432266436185
$ws.Range("F18").Value = 1.863959589624699
$ws.Range("G18").Value = 0.002504716573528116
$ws.Range("I18").Value = 1.312728397719596
$ws.Range("J18").Value = 0.1778575730178282
$ws.Range("K18").Value = 2.044013910623391
$ws.Range("L18").Value = 0.1920154464045538
$ws.Range("M18").Value = 0.5805749031424838
$ws.Range("O18").Value = 4.871276035463723
$ws.Range("C19").Value = 0.1801420462336552
$ws.Range("D19").Value = 0.1346152372537972
$ws.Range("E19").Value = 0.1360471034167485
$ws.Range("F19").Value = 1.865064322744161
$ws.Range("G19").Value = 0.002504945397868862
$ws.Range("I19").Value = 1.313557669265826
$ws.Range("J19").Value = 0.1779743004111047
$ws.Range("K19").Value = 2.030419163457964
$ws.Range("L19").Value = 0.1921075350620711
$ws.Range("M19").Value = 0.5779454781303741
$ws.Range("O19").Value = 4.875186704654226
$ws.Range("C20").Value = 0.1811899499902268
$ws.Range("D20").Value = 0.1356510401488578
$ws.Range("E20").Value = 0.1360325543227887
$ws.Range("F20").Value = 1.860164965657937
$ws.Range("G20").Value = 0.002503921967385876
$ws.Range("I20").Value = 1.309877577221222
$ws.Range("J20").Value = 0.1774533400733969
$ws.Range("K20").Value = 2.091587426991907
$ws.Range("L20").Value = 0.1916971431785512
$ws.Range("M20").Value = 0.5897812157724189
$ws.Range("O20").Value = 4.857790704504538
$ws.Range("C21").Value = 0.1847716417938301
$ws.Range("D21").Value = 0.1391632152809876
$ws.Range("E21").Value = 0.1360342640169669
$ws.Range("F21").Value = 1.844965384278368
$ws.Range("G21").Value = 0.002500593425993969
$ws.Range("I21").Value = 1.298418516835362
$ws.Range("J21").Value = 0.1757786429626478
$ws.Range("K21").Value = 2.296925125379005
$ws.Range("L21").Value = 0.1903885068482367
$ws.Range("M21").Value = 0.629598473335875
$ws.Range("O21").Value = 4.802887993591696
$ws.Range("C22").Value = 0.187157490159791
$ws.Range("D22").Value = 0.1414821538668889
$ws.Range("E22").Value = 0.1360729656793858
$ws.Range("F22").Value = 1.835972395173812
$ws.Range("G22").Value = 0.002498499715779656
$ws.Range("I22").Value = 1.291604758708878
$ws.Range("J22").Value = 0.1747404656164662
$ws.Range("K22").Value = 2.430978206119732
$ws.Range("L22").Value = 0.1895854806013446
$ws.Range("M22").Value = 0.6556564817633586
$ws.Range("O22").Value = 4.769651572835556
$ws.Range("C23").Value = 0.1858800235288527
$ws.Range("D23").Value = 0.140242361243466
$ws.Range("E23").Value = 0.1360488788004304
$ws.Range("F23").Value = 1.84068600357925
$ws.Range("G23").Value = 0.002499609715574127
$ws.Range("I23").Value = 1.295179559412148
$ws.Range("J23").Value = 0.1752894093612483
$ws.Range("K23").Value = 2.359445855281706
$ws.Range("L23").Value = 0.1900092966975748
$ws.Range("M23").Value = 0.6417457927074679
$ws.Range("O23").Value = 4.787148192171912
$ws.Range("C24").Value = 0.1811321741556071
$ws.Range("D24").Value = 0.1355940374883744
$ws.Range("E24").Value = 0.1360331626440789
$ws.Range("F24").Value = 1.860429270176056
$ws.Range("G24").Value = 0.002503977754537498
$ws.Range("I24").Value = 1.310076264839466
$ws.Range("J24").Value = 0.1774816639557244
$ws.Range("K24").Value = 2.088228998510715
$ws.Range("L24").Value = 0.1917194155090165
$ws.Range("M24").Value = 0.5891310555014968
$ws.Range("O24").Value = 4.8587326796723
$ws.Range("C25").Value = 0.1762064445356799
$ws.Range("D25").Value = 0.1306832460459191
$ws.Range("E25").Value = 0.1361783701997386
$ws.Range("F25").Value = 1.885761370770673
$ws.Range("G25").Value = 0.002509042273615065
$ws.Range("I25").Value = 1.329041377216804
$ws.Range("J25").Value = 0.180088591819576
$ws.Range("K25").Value = 1.795132476790741
$ws.Range("L25").Value = 0.1937889808269766
$ws.Range("M25").Value = 0.5325459296458774
$ws.Range("O25").Value = 4.947286404905071
